$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-04 -> 2023-10-05, i.e. 45203 -> 45204) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 443 }

$ws.Range("C2:C$lastRow").Value = 45204
